# Primer push del codigo para automatizar reportes para la DPRPOPMIYSPG
#
# Updates the "Hoja1" certificate table: several "Factor de Reajuste" /
# "Factor de Redeterminación Provisorio" values were recalculated, and the
# sheet's on-screen selection moved to H18 (scrolled toward column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Updated factor values (col J "Factor de Reajuste", plus E5 which holds a
# "Factor de Redeterminación Provisorio" outlier on row 5).
$ws.Range("J2").Value = 1.88888888
$ws.Range("J3").Value = 2.00002020202
$ws.Range("J4").Value = 3.651465416
$ws.Range("E5").Value = 2.651651
$ws.Range("J5").Value = 5.516516
$ws.Range("J6").Value = 8.65165165
$ws.Range("J7").Value = 8.651651
$ws.Range("J8").Value = 9.651651
$ws.Range("J9").Value = 5.3651651

# Move the active selection the way the author last left the sheet.
$ws.Range("H18").Select()
